$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.88", "96.881.08") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '96.881.08'
$ws.Range("E2").Value = '  -0.45%  '

# Row 3
$ws.Range("D3").Value = '3.668.94'
$ws.Range("E3").Value = '  +2.28%  '

# Row 5
$ws.Range("E5").Value = '  -0.67%  '

# Row 6
$ws.Range("D6").Value = '1.88'
$ws.Range("E6").Value = '  +11.29%  '

# Row 7
$ws.Range("D7").Value = '654.86'
$ws.Range("E7").Value = '  -0.34%  '

# Row 8
$ws.Range("D8").Value = '0.424'
$ws.Range("E8").Value = '  -0.71%  '

# Row 9
$ws.Range("E9").Value = '  +3.93%  '

# Row 10
$ws.Range("E10").Value = '  +0.04%  '

# Row 11
$ws.Range("D11").Value = '3.667.07'
$ws.Range("E11").Value = '  +2.31%  '

# Row 12
$ws.Range("D12").Value = '45.53'
$ws.Range("E12").Value = '  +2.84%  '

# Row 13
$ws.Range("E13").Value = '  +0.97%  '

# Row 14
$ws.Range("E14").Value = '  +6.26%  '

# Row 15
$ws.Range("D15").Value = '4.352.90'
$ws.Range("E15").Value = '  +2.28%  '

# Row 16
$ws.Range("D16").Value = '0.0000269'
$ws.Range("E16").Value = '  +3.27%  '

# Row 17
$ws.Range("D17").Value = '96.695.49'
$ws.Range("E17").Value = '  -0.37%  '

# Row 18
$ws.Range("D18").Value = '8.91'
$ws.Range("E18").Value = '  +3.35%  '

# Row 19
$ws.Range("D19").Value = '3.673.99'
$ws.Range("E19").Value = '  +2.25%  '

# Row 20
$ws.Range("D20").Value = '18.77'
$ws.Range("E20").Value = '  +4.15%  '

# Row 21
$ws.Range("D21").Value = '12.78'
$ws.Range("E21").Value = '  +0.43%  '

# Row 22
$ws.Range("D22").Value = '0.530'
$ws.Range("E22").Value = '  +1.05%  '

# Row 23
$ws.Range("D23").Value = '534.59'
$ws.Range("E23").Value = '  +3.70%  '

# Row 24
$ws.Range("D24").Value = '3.50'
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
$ws.Range("D25").Value = '7.21'
$ws.Range("E25").Value = '  +5.00%  '

# Row 26
$ws.Range("D26").Value = '0.0000204'
$ws.Range("E26").Value = '  -0.17%  '

# Row 27
$ws.Range("D27").Value = '102.58'
$ws.Range("E27").Value = '  +1.50%  '

# Row 28
$ws.Range("D28").Value = '13.50'
$ws.Range("E28").Value = '  +3.57%  '

# Row 29
$ws.Range("D29").Value = '3.865.95'
$ws.Range("E29").Value = '  +2.28%  '

# Row 30
$ws.Range("D30").Value = '0.167'
$ws.Range("E30").Value = '  +3.90%  '

# Row 31
$ws.Range("D31").Value = '12.39'
$ws.Range("E31").Value = '  +4.45%  '

# Row 32
$ws.Range("E32").Value = '  +1.34%  '

# Row 33
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.38%  '

# Row 34
$ws.Range("E34").Value = '  +14.32%  '

# Row 35
$ws.Range("D35").Value = '0.185'
$ws.Range("E35").Value = '  +0.89%  '

# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '32.67'
$ws.Range("E37").Value = '  +3.01%  '

# Row 38
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '658.94'
$ws.Range("E38").Value = '  +5.79%  '

# Row 39
$ws.Range("D39").Value = '0.598'
$ws.Range("E39").Value = '  +5.60%  '

# Row 40
$ws.Range("D40").Value = '8.91'
$ws.Range("E40").Value = '  +0.42%  '

# Row 41
$ws.Range("E41").Value = '  +4.60%  '

# Row 42
$ws.Range("D42").Value = '6.72'
$ws.Range("E42").Value = '  +12.21%  '

# Row 43
$ws.Range("E43").Value = '  +2.40%  '

# Row 44
$ws.Range("D44").Value = '0.962'
$ws.Range("E44").Value = '  +4.11%  '

# Row 45
$ws.Range("D45").Value = '38.75'
$ws.Range("E45").Value = '  +17.30%  '

# Row 46
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("D47").Value = '0.0459'
$ws.Range("E47").Value = '  +4.56%  '

# Row 48
$ws.Range("D48").Value = '0.445'
$ws.Range("E48").Value = '  +10.27%  '

# Row 49
$ws.Range("E49").Value = '  +1.51%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '23.62'
$ws.Range("E50").Value = '  +0.05%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '8.74'
$ws.Range("E51").Value = '  +3.06%  '
